# Applies the "begin housing benefits analysis" codebook update:
#  - adds a `label` header in column D
#  - splits the bath_shower (HH081) and toilet (HH091) yes/no questions into
#    three options each: "Yes, for sole use of the household", "Yes, shared", "No"
#  - everything below those two question blocks shifts down by two rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column D
$ws.Range("D1").Value = "label"

# Make room: the bath_shower / toilet blocks grow from 2 rows each to 3 rows
# each, so insert two new blank rows right after the current bath_shower
# block (row 35). Everything from the old row 36 onward shifts down to
# row 38 onward, which already lines up with the final layout.
$ws.Rows.Item(36).Resize(2).Insert()

# --- HH081 / bath_shower: now three response options ---
$ws.Range("A34").Value = "HH081"
$ws.Range("B34").Value = "bath_shower"
$ws.Range("C34").Value = 1
$ws.Range("D34").Value = "Yes, for sole use of the household"

$ws.Range("C35").Value = 2
$ws.Range("D35").Value = "Yes, shared"

$ws.Range("C36").Value = 3
$ws.Range("D36").Value = "No"

# --- HH091 / toilet: now three response options ---
# (the row-insert above left the old header cells sitting in row 38;
# clear them before laying out the new three-row block)
$ws.Range("A38:B38").ClearContents()

$ws.Range("A37").Value = "HH091"
$ws.Range("B37").Value = "toilet"
$ws.Range("C37").Value = 1
$ws.Range("D37").Value = "Yes, for sole use of the household"

$ws.Range("C38").Value = 2
$ws.Range("D38").Value = "Yes, shared"

$ws.Range("C39").Value = 3
$ws.Range("D39").Value = "No"

# Reflect where the author was working: scrolled/selected near the
# newly-edited rows.
$ws.Range("E38").Select()
